$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Offers" row hyperlink-display text (B4) to include the /offers path.
$ws.Range("B4").Value = "https://www.payback.in/offers"

# Move the active cell selection from B4 to B5 (matches the saved sheetView selection).
$ws.Range("B5").Select()
